# Insert a new data row right after the existing row 380 (i.e. at row 381),
# pushing the former rows 381-427 down to 382-428, and populate the new
# row 381 with the new Ají price-record data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 381; everything below shifts down by one.
$ws.Rows.Item(381).Insert()

# Fill in the new row 381 values.
$ws.Cells.Item(381, 1).Value = 9
$ws.Cells.Item(381, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(381, 3).Value = "Metropolitana"
$ws.Cells.Item(381, 4).Value = 45124
$ws.Cells.Item(381, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(381, 5).Value = 13
$ws.Cells.Item(381, 6).Value = 100112021
$ws.Cells.Item(381, 7).Value = "Ají"
$ws.Cells.Item(381, 8).Value = "Inferno"
$ws.Cells.Item(381, 9).Value = "Primera"
$ws.Cells.Item(381, 10).Value = 70
$ws.Cells.Item(381, 11).Value = 10000
$ws.Cells.Item(381, 12).Value = 17000
$ws.Cells.Item(381, 13).Value = 13500
$ws.Cells.Item(381, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(381, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(381, 16).Value = 1350
$ws.Cells.Item(381, 17).Value = 10
$ws.Cells.Item(381, 18).Value = "Hortaliza"
